$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 4 -> 5, Wrong marking -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right total 76 -> 95, Wrong total -2 -> -2.4, and the
# "obtained/max" summary text updates accordingly
$ws.Range("B12").Value = 95
$ws.Range("C12").Value = -2.4
$ws.Range("E12").Value = "92.6/140"
